# 2.1.1.1e Average per capita consumption — add the 2020 (column N) figures
# and tweak a handful of the 2019 (column M) figures to match the refreshed
# data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: thin bottom-border spacer cell under the new 2020 column ---
$ws.Range("A14").Copy()
$ws.Range("N3").PasteSpecial(-4122)

# --- Row 4: year headers ---
$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)
$ws.Cells.Item(4, 14).Value = 2020

# --- Row 5 ---
$ws.Cells.Item(5, 13).Value = 68.400000000000006
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)
$ws.Cells.Item(5, 14).Value = 68.5

# --- Row 6 ---
$ws.Cells.Item(6, 13).Value = 108.2
$ws.Range("M6").Copy()
$ws.Range("N6").PasteSpecial(-4122)
$ws.Cells.Item(6, 14).Value = 106.7

# --- Row 7 ---
$ws.Cells.Item(7, 13).Value = 51.7
$ws.Range("M7").Copy()
$ws.Range("N7").PasteSpecial(-4122)
$ws.Cells.Item(7, 14).Value = 53.2

# --- Row 8 ---
$ws.Cells.Item(8, 13).Value = 97.7
$ws.Range("M8").Copy()
$ws.Range("N8").PasteSpecial(-4122)
$ws.Cells.Item(8, 14).Value = 49.6

# --- Row 9 ---
$ws.Cells.Item(9, 12).Value = 105.6
$ws.Cells.Item(9, 13).Value = 106.7
$ws.Range("M9").Copy()
$ws.Range("N9").PasteSpecial(-4122)
$ws.Cells.Item(9, 14).Value = 108.9

# --- Row 10 ---
$ws.Cells.Item(10, 13).Value = 124.2
$ws.Range("M10").Copy()
$ws.Range("N10").PasteSpecial(-4122)
$ws.Cells.Item(10, 14).Value = 107.8

# --- Row 11 ---
$ws.Cells.Item(11, 13).Value = 138.80000000000001
$ws.Range("M11").Copy()
$ws.Range("N11").PasteSpecial(-4122)
$ws.Cells.Item(11, 14).Value = 155.69999999999999

# --- Row 12 ---
$ws.Cells.Item(12, 12).Value = 27.1
$ws.Cells.Item(12, 13).Value = 33.9
$ws.Range("M12").Copy()
$ws.Range("N12").PasteSpecial(-4122)
$ws.Cells.Item(12, 14).Value = 25.9

# --- Row 13 ---
$ws.Cells.Item(13, 13).Value = 96
$ws.Range("M13").Copy()
$ws.Range("N13").PasteSpecial(-4122)
$ws.Cells.Item(13, 14).Value = 103.5

# --- Row 14 ---
$ws.Cells.Item(14, 13).Value = 7.7
$ws.Range("M14").Copy()
$ws.Range("N14").PasteSpecial(-4122)
$ws.Cells.Item(14, 14).Value = 11

# --- Page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$excel.CutCopyMode = $false
